# books_ex.xlsx: update + add providers
# - rename the "sell price" column header to a "sell percentage" header
#   and update the two data rows accordingly (price -> percentage values)
# - tighten up a few column widths to fit the new header/values
# - move the active cell selection over to the edited column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1: "سعر البيع" (sell price) -> "نسبة البيع (%)" (sell percentage %)
$ws.Range("J1").Value = "نسبة البيع (%)"

# J2/J3: former sell-price amounts replaced with sell-percentage values
$ws.Range("J2").Value = 10
$ws.Range("J3").Value = 25

# Column width tweaks (A, E, J) to accommodate the new layout
$ws.Columns.Item(1).ColumnWidth = 55
$ws.Columns.Item(5).ColumnWidth = 10.333333333333334
$ws.Columns.Item(10).ColumnWidth = 12.666666666666666

# Active selection moves to the updated column
$ws.Range("I7").Select()
